$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1833333333333333
$ws.Range("C2").Value = 0.5333333333333333
$ws.Range("J2").Value = 0.01666666666666667
$ws.Range("P2").Value = 0.1666666666666667
$ws.Range("S2").Value = 0.1
$ws.Range("C3").Value = 0.0303030303030303
$ws.Range("J3").Value = 0.0303030303030303
$ws.Range("P3").Value = 0.8181818181818182
$ws.Range("S3").Value = 0.1212121212121212
$ws.Range("J4").Value = 0.1428571428571428
$ws.Range("P4").Value = 0.5714285714285714
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.02777777777777778
$ws.Range("F6").Value = 0.02777777777777778
$ws.Range("J6").Value = 0.3055555555555556
$ws.Range("O6").Value = 0.05555555555555555
$ws.Range("Q6").Value = 0.08333333333333333
$ws.Range("R6").Value = 0.05555555555555555
$ws.Range("S6").Value = 0.4444444444444444
$ws.Range("B7").Value = 0.2777777777777778
$ws.Range("D7").Value = 0.05555555555555555
$ws.Range("F7").Value = 0.05555555555555555
$ws.Range("J7").Value = 0.1666666666666667
$ws.Range("R7").Value = 0.05555555555555555
$ws.Range("S7").Value = 0.1666666666666667
$ws.Range("B8").Value = 0.1509433962264151
$ws.Range("F8").Value = 0.05660377358490566
$ws.Range("J8").Value = 0.05660377358490566
$ws.Range("O8").Value = 0.01886792452830189
$ws.Range("Q8").Value = 0.2641509433962264
$ws.Range("R8").Value = 0.07547169811320754
$ws.Range("S8").Value = 0.3773584905660378
$ws.Range("B9").Value = 0.103448275862069
$ws.Range("D9").Value = 0.03448275862068965
$ws.Range("F9").Value = 0.03448275862068965
$ws.Range("J9").Value = 0.1724137931034483
$ws.Range("Q9").Value = 0.3448275862068966
$ws.Range("R9").Value = 0.06896551724137931
$ws.Range("S9").Value = 0.2413793103448276
$ws.Range("B10").Value = 0.1449275362318841
$ws.Range("D10").Value = 0.02415458937198068
$ws.Range("F10").Value = 0.05797101449275362
$ws.Range("J10").Value = 0.07729468599033816
$ws.Range("O10").Value = 0.00966183574879227
$ws.Range("Q10").Value = 0.2173913043478261
$ws.Range("R10").Value = 0.1256038647342995
$ws.Range("S10").Value = 0.3429951690821256
$ws.Range("G11").Value = 0.1612903225806452
$ws.Range("K11").Value = 0.1612903225806452
$ws.Range("L11").Value = 0.6451612903225806
$ws.Range("S11").Value = 0.03225806451612903
$ws.Range("G12").Value = 0.5217391304347826
$ws.Range("J12").Value = 0.391304347826087
$ws.Range("K12").Value = 0.04347826086956522
$ws.Range("L12").Value = 0.04347826086956522
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.04761904761904762
$ws.Range("H15").Value = 0.1428571428571428
$ws.Range("I15").Value = 0.04761904761904762
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("K15").Value = 0.07142857142857142
$ws.Range("M15").Value = 0.02380952380952381
$ws.Range("O15").Value = 0.09523809523809523
$ws.Range("S15").Value = 0.2380952380952381
$ws.Range("H16").Value = 0.08333333333333333
$ws.Range("I16").Value = 0.08333333333333333
$ws.Range("J16").Value = 0.3611111111111111
$ws.Range("K16").Value = 0.1388888888888889
$ws.Range("O16").Value = 0.1388888888888889
$ws.Range("S16").Value = 0.1388888888888889
$ws.Range("F17").Value = 0.01298701298701299
$ws.Range("H17").Value = 0.1298701298701299
$ws.Range("I17").Value = 0.07792207792207792
$ws.Range("J17").Value = 0.4545454545454545
$ws.Range("K17").Value = 0.1038961038961039
$ws.Range("M17").Value = 0.01298701298701299
$ws.Range("O17").Value = 0.09090909090909091
$ws.Range("S17").Value = 0.1168831168831169
$ws.Range("H18").Value = 0.1176470588235294
$ws.Range("I18").Value = 0.08823529411764706
$ws.Range("J18").Value = 0.6176470588235294
$ws.Range("M18").Value = 0.02941176470588235
$ws.Range("O18").Value = 0.02941176470588235
$ws.Range("S18").Value = 0.1176470588235294
$ws.Range("F19").Value = 0.05747126436781609
$ws.Range("H19").Value = 0.1724137931034483
$ws.Range("I19").Value = 0.08045977011494253
$ws.Range("J19").Value = 0.4195402298850575
$ws.Range("K19").Value = 0.05747126436781609
$ws.Range("O19").Value = 0.08620689655172414
$ws.Range("S19").Value = 0.1264367816091954
